$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 1 ("Data translator"): rename the translator function, adding an
# underscore before the year -> published_SealeCarlisle_Mickes_2016()
$ws.Range("B1:M1").Value = "published_SealeCarlisle_Mickes_2016()"

# Row 8 ("maxiter"): bump the iteration cap from 20 to 2000
$ws.Range("B8:M8").Value = 2000

# Reflect the new active selection used while editing
$ws.Activate()
$ws.Range("C4:M4").Select()
